# Adding `patient.patient_problems` to Device Events (MAUDE).
#
# Inserts a new row describing the `patient.patient_problems` field right
# after the existing `Patient` section header row (old row 50 / new row 51),
# pushing every following row down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("device_event_fields")

# Insert a brand new row 51 - all rows 51.. shift down to 52..
$ws.Rows.Item(51).Insert()

$ws.Range("A51").Value() = "Patient"
$ws.Range("B51").Value() = "patient.patient_problems"
$ws.Range("C51").Value() = "array of strings"
$ws.Range("D51").Value() = "Describes actual adverse effects on the patient that may be related to the device problem observed during the reported event.`nThis is an .exact field. It has been indexed both as its exact string content, and also tokenized."

$ws.Rows.Item(51).RowHeight = 51

# View tweaks: zoom to 150% and reset the selection to A2.
$ws.Range("A2").Select()
$excel.ActiveWindow.Zoom = 150

# Page setup: scale to 29% and fit to 4x4 pages.
$ws.PageSetup.Zoom = 29
$ws.PageSetup.FitToPagesWide = 4
$ws.PageSetup.FitToPagesTall = 4
